$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 478
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 888
$ws.Range("K2").Value = 150
$ws.Range("L2").Value = 888
$ws.Range("M2").Value = -37
$ws.Range("N2").Value = -1114
$ws.Range("H9").Value = 197.2
$ws.Range("I9").Value = 245.25
$ws.Range("J9").Value = 165.16667
$ws.Range("K9").Value = 245.25
$ws.Range("L9").Value = 165.16667
$ws.Range("M9").Value = -76.25
$ws.Range("N9").Value = -503.16667
$ws.Range("H34").Value = 7577.4
$ws.Range("I34").Value = 7577.4
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 7577.4
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -7374.4
$ws.Range("H36").Value = 7577.4
$ws.Range("I36").Value = 7577.4
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 7577.4
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -6862.4
$ws.Range("H76").Value = 3013.7727
$ws.Range("I76").Value = 3017.8235
$ws.Range("K76").Value = 3017.8235
$ws.Range("M76").Value = -2702.8235
$ws.Range("H79").Value = 3013.7727
$ws.Range("I79").Value = 3017.8235
$ws.Range("K79").Value = 3017.8235
$ws.Range("M79").Value = -1925.8235
$ws.Range("H113").Value = 4098.684
$ws.Range("I113").Value = 3347.0588
$ws.Range("K113").Value = 3347.0588
$ws.Range("M113").Value = -93.05879999999979
$ws.Range("H127").Value = 111112100
$ws.Range("I127").Value = 166667200
$ws.Range("K127").Value = 500001600
$ws.Range("M127").Value = -499996640

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2614.2856
$ws.Range("I102").Value = 1933.3334
$ws.Range("J102").Value = 3125
$ws.Range("K102").Value = 1933.3334
$ws.Range("L102").Value = 3125
$ws.Range("M102").Value = -311.3334
$ws.Range("N102").Value = -6369
$ws.Range("H137").Value = 35499.832
$ws.Range("J137").Value = 35499.832
$ws.Range("L137").Value = 35499.832
$ws.Range("N137").Value = -45699.832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29169.5
$ws.Range("J55").Value = 29169.5
$ws.Range("L55").Value = 29169.5
$ws.Range("N55").Value = -29715.5
$ws.Range("H57").Value = 41600
$ws.Range("J57").Value = 41600
$ws.Range("L57").Value = 41600
$ws.Range("N57").Value = -43040
$ws.Range("H86").Value = 2245.92
$ws.Range("I86").Value = 1704.8422
$ws.Range("J86").Value = 3959.3333
$ws.Range("K86").Value = 1704.8422
$ws.Range("L86").Value = 3959.3333
$ws.Range("M86").Value = -581.8422
$ws.Range("N86").Value = -6205.3333
$ws.Range("H89").Value = 2245.92
$ws.Range("I89").Value = 1704.8422
$ws.Range("J89").Value = 3959.3333
$ws.Range("K89").Value = 8524.210999999999
$ws.Range("L89").Value = 19796.6665
$ws.Range("M89").Value = -2908.210999999999
$ws.Range("N89").Value = -31028.6665
$ws.Range("H105").Value = 2380.3062
$ws.Range("I105").Value = 2373.2896
$ws.Range("J105").Value = 2404.5454
$ws.Range("K105").Value = 2373.2896
$ws.Range("L105").Value = 2404.5454
$ws.Range("M105").Value = -626.2896000000001
$ws.Range("N105").Value = -5898.5454
$ws.Range("H134").Value = 1607034.9
$ws.Range("I134").Value = 2359404.5
$ws.Range("J134").Value = 8249.375
$ws.Range("K134").Value = 7078213.5
$ws.Range("L134").Value = 24748.125
$ws.Range("M134").Value = -7075678.5
$ws.Range("N134").Value = -29818.125
$ws.Range("H136").Value = 41600
$ws.Range("J136").Value = 41600
$ws.Range("L136").Value = 41600
$ws.Range("N136").Value = -51800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 13610.637
$ws.Range("I86").Value = 9953
$ws.Range("J86").Value = 17999.8
$ws.Range("K86").Value = 9953
$ws.Range("L86").Value = 17999.8
$ws.Range("M86").Value = -8830
$ws.Range("N86").Value = -20245.8
$ws.Range("H89").Value = 13610.637
$ws.Range("I89").Value = 9953
$ws.Range("J89").Value = 17999.8
$ws.Range("K89").Value = 49765
$ws.Range("L89").Value = 89999
$ws.Range("M89").Value = -44149
$ws.Range("N89").Value = -101231
$ws.Range("H99").Value = 2004.0646
$ws.Range("I99").Value = 1890.6
$ws.Range("J99").Value = 2210.3635
$ws.Range("K99").Value = 1890.6
$ws.Range("L99").Value = 2210.3635
$ws.Range("M99").Value = -392.5999999999999
$ws.Range("N99").Value = -5206.363499999999
$ws.Range("H126").Value = 2004.0646
$ws.Range("I126").Value = 1890.6
$ws.Range("J126").Value = 2210.3635
$ws.Range("K126").Value = 5671.799999999999
$ws.Range("L126").Value = 6631.0905
$ws.Range("M126").Value = -3201.799999999999
$ws.Range("N126").Value = -11571.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1191.8182
$ws.Range("I5").Value = 302
$ws.Range("J5").Value = 1933.3334
$ws.Range("K5").Value = 906
$ws.Range("L5").Value = 5800.0002
$ws.Range("M5").Value = -794
$ws.Range("N5").Value = -6024.0002
$ws.Range("H68").Value = 1458.2
$ws.Range("I68").Value = 1425.5
$ws.Range("J68").Value = 1480
$ws.Range("K68").Value = 4276.5
$ws.Range("L68").Value = 4440
$ws.Range("M68").Value = -3465.5
$ws.Range("N68").Value = -6062
$ws.Range("H71").Value = 1458.2
$ws.Range("I71").Value = 1425.5
$ws.Range("J71").Value = 1480
$ws.Range("K71").Value = 12829.5
$ws.Range("L71").Value = 13320
$ws.Range("M71").Value = -8773.5
$ws.Range("N71").Value = -21432
$ws.Range("H100").Value = 3532.7778
$ws.Range("J100").Value = 3532.7778
$ws.Range("L100").Value = 10598.3334
$ws.Range("N100").Value = -12220.3334
$ws.Range("H109").Value = 4680.231
$ws.Range("J109").Value = 5611.5
$ws.Range("L109").Value = 16834.5
$ws.Range("N109").Value = -18914.5
$ws.Range("H115").Value = 2834.2856
$ws.Range("J115").Value = 2834.2856
$ws.Range("L115").Value = 8502.856800000001
$ws.Range("N115").Value = -10852.8568
$ws.Range("H135").Value = 1191.8182
$ws.Range("I135").Value = 302
$ws.Range("J135").Value = 1933.3334
$ws.Range("K135").Value = 2718
$ws.Range("L135").Value = 17400.0006
$ws.Range("M135").Value = -183
$ws.Range("N135").Value = -22470.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 29000
$ws.Range("J51").Value = 29000
$ws.Range("L51").Value = 29000
$ws.Range("N51").Value = -30018
$ws.Range("H57").Value = 18819.8
$ws.Range("J57").Value = 18819.8
$ws.Range("L57").Value = 18819.8
$ws.Range("N57").Value = -20459.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2754.6
$ws.Range("I61").Value = 2189.0588
$ws.Range("J61").Value = 3288.7222
$ws.Range("K61").Value = 2189.0588
$ws.Range("L61").Value = 3288.7222
$ws.Range("M61").Value = -1987.0588
$ws.Range("N61").Value = -3692.7222
$ws.Range("H64").Value = 30429.6
$ws.Range("J64").Value = 30429.6
$ws.Range("L64").Value = 30429.6
$ws.Range("N64").Value = -30879.6
$ws.Range("H67").Value = 30429.6
$ws.Range("J67").Value = 30429.6
$ws.Range("L67").Value = 30429.6
$ws.Range("N67").Value = -31989.6
$ws.Range("H113").Value = 2754.6
$ws.Range("I113").Value = 2189.0588
$ws.Range("J113").Value = 3288.7222
$ws.Range("K113").Value = 2189.0588
$ws.Range("L113").Value = 3288.7222
$ws.Range("M113").Value = -19.05879999999979
$ws.Range("N113").Value = -7628.7222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 101813.25
$ws.Range("I122").Value = 1843.4286
$ws.Range("J122").Value = 241771
$ws.Range("K122").Value = 5530.2858
$ws.Range("L122").Value = 725313
$ws.Range("M122").Value = -3080.2858
$ws.Range("N122").Value = -730213
